$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the raw input values for the "L3" row (row 6)
$ws.Range("C6").Value = 1180
$ws.Range("F6").Value = 1700
$ws.Range("G6").Value = 2150

# Move the active selection to H15 (matches the saved sheet view state)
$ws.Range("H15").Select()
